# RPAR_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclosure footer (A18)
#  - refresh the Weight (D) / Percent Change (E) figures for rows 2-15
#
# The worksheet is protected, so it must be unprotected before writing and
# re-protected afterwards to restore the original "protected" state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# --- Footer disclosure text: 2021-05-10 -> 2021-05-11 -------------------
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."
# Writing a 2-line value bumps the row's auto height; put it back so the
# row element doesn't pick up a spurious ht/customHeight attribute.
$ws.Rows.Item(18).AutoFit()

# --- Weight (D) / Percent Change (E) refresh -----------------------------
$ws.Range("D2").Value = 0.05667975857747388
$ws.Range("E2").Value = -0.008048848182070523

$ws.Range("D3").Value = 0.02379944206436186
$ws.Range("E3").Value = -0.0119529593213803

$ws.Range("D4").Value = 0.0309709338605056
$ws.Range("E4").Value = -0.002666666666666595

$ws.Range("D5").Value = 0.03317724419647305
$ws.Range("E5").Value = -0.02327327327327333

$ws.Range("D6").Value = 0.03998101760328888
$ws.Range("E6").Value = 0.002951013181192064

$ws.Range("D7").Value = 0.01953692390785672
$ws.Range("E7").Value = -0.009159229685298365

$ws.Range("D8").Value = 0.004204181371032443
$ws.Range("E8").Value = 0.002338634237605053

$ws.Range("D9").Value = 0.006992548622727488
$ws.Range("E9").Value = -0.008998875140607487

$ws.Range("D10").Value = 0.07181323838639142
$ws.Range("E10").Value = 0

$ws.Range("D11").Value = 0.0718918947263546
$ws.Range("E11").Value = 0.0005470459518597259

$ws.Range("D12").Value = 0.1445388903163558
$ws.Range("E12").Value = -0.005659555942533756

$ws.Range("D13").Value = 0.3819591196782431
$ws.Range("E13").Value = -0.001662874146683091

$ws.Range("D14").Value = 0.1144548066889352
$ws.Range("E14").Value = -0.00541190619362597

$ws.Range("E15").Value = -0.003742730843248387

$ws.Protect()
